$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 240
$ws1.Range("F6").Value = 81
$ws1.Range("F7").Value = 810
$ws1.Range("F8").Value = 460
$ws1.Range("F9").Value = 66
$ws1.Range("F12").Value = 215
$ws1.Range("F19").Value = 7380
$ws1.Range("F22").Value = 3344
$ws1.Range("F23").Value = 773
$ws1.Range("F29").Value = 1413
$ws1.Range("F31").Value = 47
$ws1.Range("F33").Value = 1092
$ws1.Range("F34").Value = 1524

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 74

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 240
$ws4.Range("F9").Value = 81
$ws4.Range("F10").Value = 810
$ws4.Range("F11").Value = 460
$ws4.Range("F12").Value = 66
$ws4.Range("F16").Value = 215
$ws4.Range("F23").Value = 7380
$ws4.Range("F26").Value = 3344
$ws4.Range("F27").Value = 773
$ws4.Range("F34").Value = 1413
$ws4.Range("F36").Value = 47
$ws4.Range("F38").Value = 1092
$ws4.Range("F39").Value = 1524
$ws4.Range("F43").Value = 74
